$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.759.68'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '2.460.34'
$ws.Range("E3").Value = '  +0.81%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = "'573.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.67%  '

$ws.Range("D6").Value = "'145.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '2.458.96'
$ws.Range("E9").Value = '  +0.85%  '

$ws.Range("E10").Value = '  +1.04%  '

$ws.Range("D11").Value = "'0.163"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.69%  '

$ws.Range("D12").Value = "'5.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.60%  '

$ws.Range("E13").Value = '  +0.77%  '

$ws.Range("D14").Value = "'28.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("E15").Value = '  -0.51%  '

$ws.Range("D16").Value = '2.906.73'
$ws.Range("E16").Value = '  +0.81%  '

$ws.Range("D17").Value = '62.693.22'
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").Value = '2.450.04'
$ws.Range("E18").Value = '  +0.46%  '

$ws.Range("D19").Value = "'7.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.56%  '

$ws.Range("E20").Value = '  +0.53%  '

$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").Value = "'4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").Value = "'2.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.74%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("D25").Value = "'10.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +18.47%  '

$ws.Range("D26").Value = "'65.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("D27").Value = "'652.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("E28").Value = '  +0.97%  '

$ws.Range("D29").Value = '0.0₃0975'
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -12.12%  '

$ws.Range("E31").Value = '  +2.97%  '

$ws.Range("D32").Value = "'7.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.51%  '

$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("E34").Value = '  -3.90%  '

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("E36").Value = '  +2.54%  '

$ws.Range("E37").Value = '  +0.32%  '

$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = "'0.369"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.82%  '

$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = "'18.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.63%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = "'150.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.22%  '

$ws.Range("E41").Value = '  -2.00%  '

$ws.Range("D42").Value = "'2.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.85%  '

$ws.Range("E43").Value = '  -1.51%  '

$ws.Range("D44").Value = '0.0₆0314'
$ws.Range("E44").Value = '  -57.30%  '

$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").Value = "'153.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.75%  '

$ws.Range("E47").Value = '  +1.35%  '

$ws.Range("E48").Value = '  -0.10%  '

$ws.Range("D49").Value = "'0.606"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.74%  '

$ws.Range("D50").Value = "'20.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.41%  '

$ws.Range("D51").Value = "'0.0509"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.35%  '
